$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed date) column C for rows 2-8 from 45184 (2023-09-15)
# to 45185 (2023-09-16), keeping the existing date format.
$ws.Range("C2:C8").Value = 45185
